$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 74-121 with revised prices ---
$ws.Range("B74").Value = 251.14
$ws.Range("B75").Value = 239.41
$ws.Range("B76").Value = 251.67
$ws.Range("B77").Value = 229.36
$ws.Range("B78").Value = 227.96
$ws.Range("B79").Value = 232.58
$ws.Range("B80").Value = 230.68
$ws.Range("B81").Value = 243.53
$ws.Range("B82").Value = 277.92
$ws.Range("B83").Value = 306.87
$ws.Range("B84").Value = 319.53
$ws.Range("B85").Value = 325.72
$ws.Range("B86").Value = 323.89
$ws.Range("B87").Value = 320.36
$ws.Range("B88").Value = 322.86
$ws.Range("B89").Value = 324.82
$ws.Range("B90").Value = 342.72
$ws.Range("B91").Value = 366.48
$ws.Range("B92").Value = 321.58
$ws.Range("B93").Value = 318.57
$ws.Range("B94").Value = 284.2
$ws.Range("B95").Value = 293.88
$ws.Range("B96").Value = 280.48
$ws.Range("B97").Value = 261.43
$ws.Range("B98").Value = 223.34
$ws.Range("B99").Value = 208.13
$ws.Range("B100").Value = 195.94
$ws.Range("B101").Value = 193.03
$ws.Range("B102").Value = 178.16
$ws.Range("B103").Value = 266.03
$ws.Range("B104").Value = 287.04
$ws.Range("B105").Value = 323.39
$ws.Range("B106").Value = 314.19
$ws.Range("B107").Value = 326.59
$ws.Range("B108").Value = 340.76
$ws.Range("B109").Value = 358.04
$ws.Range("B110").Value = 352.38
$ws.Range("B111").Value = 348.35
$ws.Range("B112").Value = 365.69
$ws.Range("B113").Value = 446.42
$ws.Range("B114").Value = 457.65
$ws.Range("B115").Value = 484.41
$ws.Range("B116").Value = 460.22
$ws.Range("B117").Value = 454.66
$ws.Range("B118").Value = 416.82
$ws.Range("B119").Value = 329.78
$ws.Range("B120").Value = 313.89
$ws.Range("B121").Value = 299.3

# --- Append new rows 122-169 for Dec 6-7, 2022 ---
$ws.Range("A122").Value = "2022-12-06 00:00"
$ws.Range("B122").Value = 208.46
$ws.Range("A123").Value = "2022-12-06 01:00"
$ws.Range("B123").Value = 152.32
$ws.Range("A124").Value = "2022-12-06 02:00"
$ws.Range("B124").Value = 96.87
$ws.Range("A125").Value = "2022-12-06 03:00"
$ws.Range("B125").Value = 98.98
$ws.Range("A126").Value = "2022-12-06 04:00"
$ws.Range("B126").Value = 152.48
$ws.Range("A127").Value = "2022-12-06 05:00"
$ws.Range("B127").Value = 300.26
$ws.Range("A128").Value = "2022-12-06 06:00"
$ws.Range("B128").Value = 325.48
$ws.Range("A129").Value = "2022-12-06 07:00"
$ws.Range("B129").Value = 491.85
$ws.Range("A130").Value = "2022-12-06 08:00"
$ws.Range("B130").Value = 484.14
$ws.Range("A131").Value = "2022-12-06 09:00"
$ws.Range("B131").Value = 488.66
$ws.Range("A132").Value = "2022-12-06 10:00"
$ws.Range("B132").Value = 493.34
$ws.Range("A133").Value = "2022-12-06 11:00"
$ws.Range("B133").Value = 490.45
$ws.Range("A134").Value = "2022-12-06 12:00"
$ws.Range("B134").Value = 494.17
$ws.Range("A135").Value = "2022-12-06 13:00"
$ws.Range("B135").Value = 485.61
$ws.Range("A136").Value = "2022-12-06 14:00"
$ws.Range("B136").Value = 499.27
$ws.Range("A137").Value = "2022-12-06 15:00"
$ws.Range("B137").Value = 494.03
$ws.Range("A138").Value = "2022-12-06 16:00"
$ws.Range("B138").Value = 502.4
$ws.Range("A139").Value = "2022-12-06 17:00"
$ws.Range("B139").Value = 521.86
$ws.Range("A140").Value = "2022-12-06 18:00"
$ws.Range("B140").Value = 494.13
$ws.Range("A141").Value = "2022-12-06 19:00"
$ws.Range("B141").Value = 455.3
$ws.Range("A142").Value = "2022-12-06 20:00"
$ws.Range("B142").Value = 401.62
$ws.Range("A143").Value = "2022-12-06 21:00"
$ws.Range("B143").Value = 347.25
$ws.Range("A144").Value = "2022-12-06 22:00"
$ws.Range("B144").Value = 326.56
$ws.Range("A145").Value = "2022-12-06 23:00"
$ws.Range("B145").Value = 290.8
$ws.Range("A146").Value = "2022-12-07 00:00"
$ws.Range("B146").Value = 288.2
$ws.Range("A147").Value = "2022-12-07 01:00"
$ws.Range("B147").Value = 256.76
$ws.Range("A148").Value = "2022-12-07 02:00"
$ws.Range("B148").Value = 222.11
$ws.Range("A149").Value = "2022-12-07 03:00"
$ws.Range("B149").Value = 217.81
$ws.Range("A150").Value = "2022-12-07 04:00"
$ws.Range("B150").Value = 250.72
$ws.Range("A151").Value = "2022-12-07 05:00"
$ws.Range("B151").Value = 282.55
$ws.Range("A152").Value = "2022-12-07 06:00"
$ws.Range("B152").Value = 307.08
$ws.Range("A153").Value = "2022-12-07 07:00"
$ws.Range("B153").Value = 413.82
$ws.Range("A154").Value = "2022-12-07 08:00"
$ws.Range("B154").Value = 441.29
$ws.Range("A155").Value = "2022-12-07 09:00"
$ws.Range("B155").Value = 429.48
$ws.Range("A156").Value = "2022-12-07 10:00"
$ws.Range("B156").Value = 412.26
$ws.Range("A157").Value = "2022-12-07 11:00"
$ws.Range("B157").Value = 387.29
$ws.Range("A158").Value = "2022-12-07 12:00"
$ws.Range("B158").Value = 413.89
$ws.Range("A159").Value = "2022-12-07 13:00"
$ws.Range("B159").Value = 374.49
$ws.Range("A160").Value = "2022-12-07 14:00"
$ws.Range("B160").Value = 444.48
$ws.Range("A161").Value = "2022-12-07 15:00"
$ws.Range("B161").Value = 475.72
$ws.Range("A162").Value = "2022-12-07 16:00"
$ws.Range("B162").Value = 494.59
$ws.Range("A163").Value = "2022-12-07 17:00"
$ws.Range("B163").Value = 456.68
$ws.Range("A164").Value = "2022-12-07 18:00"
$ws.Range("B164").Value = 435.7
$ws.Range("A165").Value = "2022-12-07 19:00"
$ws.Range("B165").Value = 423.04
$ws.Range("A166").Value = "2022-12-07 20:00"
$ws.Range("B166").Value = 445.31
$ws.Range("A167").Value = "2022-12-07 21:00"
$ws.Range("B167").Value = 356.35
$ws.Range("A168").Value = "2022-12-07 22:00"
$ws.Range("B168").Value = 319.7
$ws.Range("A169").Value = "2022-12-07 23:00"
$ws.Range("B169").Value = 296.73

# --- Resize the table to cover the new range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B170"))
